$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Resveratrol (Simulated)" data row (row 10)
$ws.Range("A10").Value = "Resveratrol (Simulated)"
$ws.Range("B10").Value = 19.483
$ws.Range("C10").Value = 8.553
$ws.Range("D10").Value = 15.297
$ws.Range("E10").Value = "Simulated - CosmoQuick QSPR"

# E1 header cell now shares the same (non-wrapped-border) look as the other
# header cells B1:D1 instead of its own distinct bordered/wrapped style
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Give every data cell (rows 2-10, columns A-E) a thin box border
$rng = $ws.Range("A2:E10")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Restore selection like the saved workbook
$ws.Range("E30").Select()
